$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" header column in H, matching the style of the existing
# header cells (e.g. G1) by copying its formatting, then setting the text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data values for H2:H6 (plain numbers, default/no special style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
